$d = $word.ActiveDocument

# The <id>...</id> tag for each page was split across three runs (the
# literal "<id>"/"</id>" runs in the Courier New / dark-yellow "tag"
# style, with the bare id value "p036vN_n" run sandwiched between them in
# plain black text). Re-download normalized each of these onto a single
# run. Search for the full "<id>p036v_n</id>" text (which already spans
# all three runs) and replace it with itself: Word's Find/Replace merges
# the matched range into one run, taking on the formatting of the first
# run in the match (the Courier New "tag" formatting), which is exactly
# the desired end state.
foreach ($n in 1..3) {
    $tag = "<id>p036v_$n</id>"
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($tag, $true, $false, $false, $false, $false, $true, 1, $false, $tag, 2)
}
